# Generate Report for Handoff
# Update "Latest Handoff Date/Datetime" timestamps for the rows that are
# re-handed-off (3adb1181-227b-4842-8e86-e1784b95ff51 and the 5 "Ready for
# handoff" rows that shared its previous handoff timestamp).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $overviewRows) {
    $wsOverview.Range("D$r").Value = "2016-03-18 17:37:02"
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $zhCnRows) {
    $wsZhCn.Range("E$r").Value = "2016-03-18 17:36:54"
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $deDeRows) {
    $wsDeDe.Range("E$r").Value = "2016-03-18 17:37:02"
}
